$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01846133333333333
$ws.Range("H2").Value = 0.055384
$ws.Range("I2").Value = 0.005172740524168673
$ws.Range("J2").Value = 0.005172740524168674
$ws.Range("M2").Value = 8.488196666666667
$ws.Range("N2").Value = 25.46459
$ws.Range("O2").Value = 0.04138402976425696
$ws.Range("P2").Value = 0.04138402976425696
$ws.Range("Q2").Value = 0.1567034280622222
$ws.Range("R2").Value = 1.41033085256
$ws.Range("S2").Value = 0.0002140688478149745
$ws.Range("T2").Value = 0.0002140688478149745
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01846133333333333
$ws.Range("H3").Value = 0.055384
$ws.Range("I3").Value = 0.005172740524168673
$ws.Range("J3").Value = 0.005172740524168674
$ws.Range("O3").Value = 0.3297460182766552
$ws.Range("P3").Value = 0.3297460182766552
$ws.Range("Q3").Value = 1.248605603373333
$ws.Range("R3").Value = 11.23745043036
$ws.Range("S3").Value = 0.001705690591422918
$ws.Range("T3").Value = 0.001705690591422918
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01846133333333333
$ws.Range("H4").Value = 0.055384
$ws.Range("I4").Value = 0.005172740524168673
$ws.Range("J4").Value = 0.005172740524168674
$ws.Range("O4").Value = 0.6288699519590879
$ws.Range("P4").Value = 0.6288699519590879
$ws.Range("Q4").Value = 2.381258612046222
$ws.Range("R4").Value = 21.431327508416
$ws.Range("S4").Value = 0.003252981084930781
$ws.Range("T4").Value = 0.003252981084930781
$ws.Range("I5").Value = 0.5495916323842472
$ws.Range("J5").Value = 0.5495916323842472
$ws.Range("M5").Value = 8.488196666666667
$ws.Range("N5").Value = 25.46459
$ws.Range("O5").Value = 0.04138402976425696
$ws.Range("P5").Value = 0.04138402976425696
$ws.Range("Q5").Value = 16.64937423915444
$ws.Range("R5").Value = 149.84436815239
$ws.Range("S5").Value = 0.02274431647277626
$ws.Range("T5").Value = 0.02274431647277626
$ws.Range("I6").Value = 0.5495916323842472
$ws.Range("J6").Value = 0.5495916323842472
$ws.Range("O6").Value = 0.3297460182766552
$ws.Range("P6").Value = 0.3297460182766552
$ws.Range("S6").Value = 0.1812256524568727
$ws.Range("T6").Value = 0.1812256524568727
$ws.Range("I7").Value = 0.5495916323842472
$ws.Range("J7").Value = 0.5495916323842472
$ws.Range("O7").Value = 0.6288699519590879
$ws.Range("P7").Value = 0.6288699519590879
$ws.Range("S7").Value = 0.3456216634545983
$ws.Range("T7").Value = 0.3456216634545983
$ws.Range("I8").Value = 0.445235627091584
$ws.Range("J8").Value = 0.445235627091584
$ws.Range("M8").Value = 8.488196666666667
$ws.Range("N8").Value = 25.46459
$ws.Range("O8").Value = 0.04138402976425696
$ws.Range("P8").Value = 0.04138402976425696
$ws.Range("Q8").Value = 13.48800480803111
$ws.Range("R8").Value = 121.39204327228
$ws.Range("S8").Value = 0.01842564444366573
$ws.Range("T8").Value = 0.01842564444366573
$ws.Range("I9").Value = 0.445235627091584
$ws.Range("J9").Value = 0.445235627091584
$ws.Range("O9").Value = 0.3297460182766552
$ws.Range("P9").Value = 0.3297460182766552
$ws.Range("S9").Value = 0.1468146752283595
$ws.Range("T9").Value = 0.1468146752283595
$ws.Range("I10").Value = 0.445235627091584
$ws.Range("J10").Value = 0.445235627091584
$ws.Range("O10").Value = 0.6288699519590879
$ws.Range("P10").Value = 0.6288699519590879
$ws.Range("S10").Value = 0.2799953074195589
$ws.Range("T10").Value = 0.2799953074195589
